$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the existing "JQW, JQW1" Qwiic connector row (row 11).
# This shifts the old rows 12-21 down to 13-22 and copies row 11's formatting.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the corrected Qwiic JST connector part.
$ws.Range("A12").Value = "Sparkfun #14417"
$ws.Range("B12").Value = "JQW, JQW1"
$ws.Range("C12").Value = "Qwiic JST Connector Horizontal"
$ws.Range("D12").Value = "JST04_1MM_RA"
$ws.Range("E12").Value = "Qwiic Right Angle"
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 2

# Mark the old (wrong) part in row 11 as no longer valid by striking it through.
$ws.Range("A11:G11").Font.Strikethrough = $true

# Grow the query table / autofilter range to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G22"))

# Keep the external-data defined name range in sync with the table.
$wb.Names.Item(1).RefersTo = "=T41SimpleLSSArduinoShieldV02!`$B`$1:`$G`$22"

# Match the author's last-saved selection.
$ws.Range("C9").Select()
